$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp banner (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 21 de Mayo de 2020 a las 16:35"

# --- Row 4: Estados Unidos - updated counts ---
$ws.Range("B4").Value = 1596526
$ws.Range("C4").Value = 3803
$ws.Range("E4").Value = 1130496
$ws.Range("G4").Value = 121
$ws.Range("H4").Value = 95057

# --- Row 54: Noruega - updated counts ---
$ws.Range("B54").Value = 8301
$ws.Range("C54").Value = 20
$ws.Range("E54").Value = 8034
$ws.Range("G54").Value = 1
$ws.Range("H54").Value = 235

# --- Row 76: Uzbekistan - updated counts ---
$ws.Range("B76").Value = 2964
$ws.Range("C76").Value = 25
$ws.Range("E76").Value = 544

# --- Rows 80-85: Tayikistan moves up (re-sorted by Casos totales) pushing
#     Bosnia y Herzegovina, Bulgaria, Guatemala, Croacia and Costa de Marfil
#     each one row down. Row 79 (Senegal) and row 86 (Cuba) stay put.

# Row 80 becomes Tayikistan with freshly updated figures
$ws.Range("A80").Value = "Tayikistan"
$ws.Range("B80").Value = 2350
$ws.Range("C80").Value = 210
$ws.Range("D80").Value = 470
$ws.Range("E80").Value = 1836
$ws.Range("F80").Value = 0
$ws.Range("G80").Value = 3
$ws.Range("H80").Value = 44

# Row 81 becomes Bosnia y Herzegovina (former row 80 data)
$ws.Range("A81").Value = "Bosnia y Herzegovina"
$ws.Range("B81").Value = 2350
$ws.Range("C81").Value = 12
$ws.Range("D81").Value = 1596
$ws.Range("E81").Value = 614
$ws.Range("F81").Value = 0
$ws.Range("G81").Value = 4
$ws.Range("H81").Value = 140

# Row 82 becomes Bulgaria (former row 81 data)
$ws.Range("A82").Value = "Bulgaria"
$ws.Range("B82").Value = 2331
$ws.Range("C82").Value = 39
$ws.Range("D82").Value = 727
$ws.Range("E82").Value = 1484
$ws.Range("F82").Value = 0
$ws.Range("G82").Value = 4
$ws.Range("H82").Value = 120

# Row 83 becomes Guatemala (former row 82 data)
$ws.Range("A83").Value = "Guatemala"
$ws.Range("B83").Value = 2265
$ws.Range("C83").Value = 132
$ws.Range("D83").Value = 159
$ws.Range("E83").Value = 2061
$ws.Range("F83").Value = 0
$ws.Range("G83").Value = 2
$ws.Range("H83").Value = 45

# Row 84 becomes Croacia (former row 83 data)
$ws.Range("A84").Value = "Croacia"
$ws.Range("B84").Value = 2237
$ws.Range("C84").Value = 3
$ws.Range("D84").Value = 1978
$ws.Range("E84").Value = 162
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 1
$ws.Range("H84").Value = 97

# Row 85 becomes Costa de Marfil (former row 84 data)
$ws.Range("A85").Value = "Costa de Marfil"
$ws.Range("B85").Value = 2231
$ws.Range("C85").Value = 0
$ws.Range("D85").Value = 1083
$ws.Range("E85").Value = 1119
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 29
